$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing data range B2:B5
$ws.Range("B2:B5").ClearContents()

# Set the new values
$ws.Range("B1").Value = "Informação1"
$ws.Range("B3").Value = "Informação3"
